# Generate Report for Handback
# Rewrites the per-run identifiers (file UUIDs, content hashes, timestamps)
# that appear across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$UUID1_OLD = "046c4de4-13fc-4e10-9864-81ba3df330e1"
$UUID1_NEW = "3e4b9cfa-e88c-4031-80e9-09febb2e63d4"
$UUID2_OLD = "556897de-3e1b-4377-bad9-7449368134b2"
$UUID2_NEW = "ffff532e7da8-8e06-4a8e-8b4e-7bd5bc186051"

$HASH_NEW = "a6c83b066a3ae09bb2afd8db03a9ed7680c3d691"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$UUID1_NEW.md"
$wsOverview.Range("B2").Value = "e2e\$UUID1_NEW.md"
$wsOverview.Range("G2").Value = "2017-01-03 08:07:20"

$wsOverview.Range("A3").Value = "$UUID2_NEW.md"
$wsOverview.Range("B3").Value = "e2e\$UUID2_NEW.md"
$wsOverview.Range("G3").Value = "2017-01-03 08:07:20"

foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = "e2e\$UUID1_NEW.md"
    } elseif ($addr -eq '$B$3') {
        $hl.TextToDisplay = "e2e\$UUID2_NEW.md"
    }
}

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$UUID1_NEW.md"
$wsZhCn.Range("G2").Value = "$UUID1_NEW.$HASH_NEW.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2017-01-03 08:07:09"
$wsZhCn.Range("J2").Value = "$UUID1_NEW.md"
$wsZhCn.Range("K2").Value = "$UUID1_NEW.$HASH_NEW.zh-cn.xlf"
$wsZhCn.Range("L2").Value = "2017-01-03 08:07:44"

$wsZhCn.Range("A3").Value = "$UUID2_NEW.md"
$wsZhCn.Range("G3").Value = "$UUID1_NEW.$HASH_NEW.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2017-01-03 08:07:09"
$wsZhCn.Range("J3").Value = "$UUID2_NEW.md"
$wsZhCn.Range("K3").Value = "$UUID1_NEW.$HASH_NEW.zh-cn.xlf"
$wsZhCn.Range("L3").Value = "2017-01-03 08:07:44"

foreach ($hl in $wsZhCn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "$UUID1_NEW.md"
    } elseif ($addr -eq '$J$2') {
        $hl.TextToDisplay = "$UUID1_NEW.md"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "$UUID2_NEW.md"
    } elseif ($addr -eq '$J$3') {
        $hl.TextToDisplay = "$UUID2_NEW.md"
    }
}

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$UUID1_NEW.md"
$wsDeDe.Range("G2").Value = "$UUID1_NEW.$HASH_NEW.de-de.xlf"
$wsDeDe.Range("H2").Value = "2017-01-03 08:07:20"
$wsDeDe.Range("J2").Value = "$UUID1_NEW.md"
$wsDeDe.Range("K2").Value = "$UUID1_NEW.$HASH_NEW.de-de.xlf"
$wsDeDe.Range("L2").Value = "2017-01-03 08:07:57"

$wsDeDe.Range("A3").Value = "$UUID2_NEW.md"
$wsDeDe.Range("G3").Value = "$UUID1_NEW.$HASH_NEW.de-de.xlf"
$wsDeDe.Range("H3").Value = "2017-01-03 08:07:20"
$wsDeDe.Range("J3").Value = "$UUID2_NEW.md"
$wsDeDe.Range("K3").Value = "$UUID1_NEW.$HASH_NEW.de-de.xlf"
$wsDeDe.Range("L3").Value = "2017-01-03 08:07:57"

foreach ($hl in $wsDeDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "$UUID1_NEW.md"
    } elseif ($addr -eq '$J$2') {
        $hl.TextToDisplay = "$UUID1_NEW.md"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "$UUID2_NEW.md"
    } elseif ($addr -eq '$J$3') {
        $hl.TextToDisplay = "$UUID2_NEW.md"
    }
}
